# Patches applied as on 24th Nov
# Appends rows 11-13 to Sheet1's patch log, matching the date-formatted
# "B" column styles already used by rows 2-10, then updates the view
# (selection / scroll) to reflect where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11: 15-Nov-2011, EKA_METALS_PATCH_0214.txt -----------------------
$ws.Range("B11").NumberFormat = "d-mmm-yy"
$ws.Range("B11").Value = 40862
$ws.Range("I11").Value = "EKA_METALS_PATCH_0214.txt"

# --- Row 12: 23-Nov-2011, CDC 1.3.3 / 1.4.9 patch --------------------------
$ws.Range("B12").NumberFormat = "d-mmm"
$ws.Range("B12").Value = 40870
$ws.Range("D12").Value = "1.4.9"
$ws.Range("E12").Value = "CDC 1.3.3"
$ws.Range("F12").Value = "CDC 1.3.3"
$ws.Range("I12").Value = "EKA_METALS_PATCH_0234.txt"

# --- Row 13: 24-Nov-2011, EKA_METALS_PATCH_0237.txt ------------------------
$ws.Range("B13").NumberFormat = "d-mmm"
$ws.Range("B13").Value = 40871
$ws.Range("I13").Value = "EKA_METALS_PATCH_0237.txt"

# Column B widened slightly to keep fitting the new dates.
$ws.Columns.Item(2).ColumnWidth = 9 + 1/6

# Leave the view scrolled/selected on the last entry added.
$ws.Range("I13").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
